$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, re-using the same formatting as the other
# header cells (bold font, thin border, centered/top aligned) by copying
# the format from G1 rather than re-building it (keeps the same style
# index instead of creating a near-duplicate style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column - all rows are 1 for this sheet.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
